# job-and-interview.xlsx — "import job interview excel"
#
# The sheet used to carry separate "vendor_first_name" / "vendor_last_name"
# columns (J:K, with sample values "New" / "Vendor"). This import drops that
# split-name pair entirely (the vendor is now only tracked by e-mail), fixes
# a header typo ("benefits" -> "benifits"), and turns the vendor e-mail cell
# into a live mailto: hyperlink, as Excel does automatically when an e-mail
# address is entered/confirmed in a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the vendor_first_name / vendor_last_name columns (J:K) entirely,
# shifting everything to their right (vendor_email, position,
# interview_start_date, interview_end_date) two columns to the left.
$ws.Range("J1:K2").Delete(-4159)  # xlShiftToLeft

# Header typo fix: "benefits" -> "benifits"
$ws.Range("D1").Value = "benifits"

# The vendor e-mail address (now in J2, after the column shift) becomes a
# clickable mailto: hyperlink, picking up Excel's built-in "Hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:vendor@yopmai.com", "", "", "vendor@yopmai.com")

# Leave the selection where the editor last left it when saving.
$ws.Range("E7").Select() | Out-Null
